$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 and 14 swap entirely (Litecoin / WrappedEther swap positions)
$ws.Range("B13").Value = '''Litecoin'
$ws.Range("C13").Value = '''https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D13").Value = '''101.69'
$ws.Range("E13").Value = '''  +1.67%  '

$ws.Range("B14").Value = '''WrappedEther'
$ws.Range("C14").Value = '''https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '''2.003.33'
$ws.Range("E14").Value = '''  +7.15%  '

# Update Price (D) and Volume(1h) (E) columns for remaining rows
$ws.Range("D2").Value = '''31.298.65'
$ws.Range("E2").Value = '''  +3.14%  '
$ws.Range("D3").Value = '''2.004.49'
$ws.Range("E3").Value = '''  +7.23%  '
$ws.Range("D4").Value = '''0.9997'
$ws.Range("E4").Value = '''  -0.01%  '
$ws.Range("D5").Value = '''0.7931'
$ws.Range("E5").Value = '''  +68.04%  '
$ws.Range("D6").Value = '''258.64'
$ws.Range("E6").Value = '''  +5.97%  '
$ws.Range("D7").Value = '''0.9988'
$ws.Range("E7").Value = '''  -0.11%  '
$ws.Range("D8").Value = '''0.3623'
$ws.Range("E8").Value = '''  +25.97%  '
$ws.Range("D9").Value = '''28.45'
$ws.Range("E9").Value = '''  +30.40%  '
$ws.Range("D10").Value = '''0.07085'
$ws.Range("E10").Value = '''  +9.25%  '
$ws.Range("D11").Value = '''0.8546'
$ws.Range("E11").Value = '''  +17.54%  '
$ws.Range("D12").Value = '''0.08181'
$ws.Range("E12").Value = '''  +4.96%  '
$ws.Range("D15").Value = '''5.617'
$ws.Range("E15").Value = '''  +8.71%  '
$ws.Range("D16").Value = '''275.86'
$ws.Range("E16").Value = '''  -2.60%  '
$ws.Range("D17").Value = '''14.98'
$ws.Range("E17").Value = '''  +14.39%  '
$ws.Range("D18").Value = '''31.296.49'
$ws.Range("E18").Value = '''  +3.20%  '
$ws.Range("D19").Value = '''5.935'
$ws.Range("E19").Value = '''  +12.49%  '
$ws.Range("D20").Value = '''0.000007989'
$ws.Range("E20").Value = '''  +6.79%  '
$ws.Range("D21").Value = '''2.268.34'
$ws.Range("E21").Value = '''  +7.48%  '
$ws.Range("D22").Value = '''0.9988'
$ws.Range("E22").Value = '''  -0.13%  '
$ws.Range("D23").Value = '''0.9998'
$ws.Range("E23").Value = '''  +0.01%  '
$ws.Range("D24").Value = '''7.193'
$ws.Range("E24").Value = '''  +14.41%  '
$ws.Range("D25").Value = '''10.12'
$ws.Range("E25").Value = '''  +11.95%  '
$ws.Range("D26").Value = '''0.1522'
$ws.Range("E26").Value = '''  +57.84%  '
$ws.Range("D27").Value = '''164.83'
$ws.Range("E27").Value = '''  +1.06%  '
$ws.Range("D28").Value = '''20.08'
$ws.Range("E28").Value = '''  +5.63%  '
$ws.Range("E29").Value = '''  +26.84%  '
$ws.Range("D30").Value = '''1.620'
$ws.Range("E30").Value = '''  +9.30%  '
$ws.Range("D31").Value = '''4.624'
$ws.Range("E31").Value = '''  +9.31%  '
$ws.Range("D32").Value = '''1.358'
$ws.Range("E32").Value = '''  +2.95%  '
$ws.Range("D33").Value = '''4.428'
$ws.Range("E33").Value = '''  +6.97%  '
$ws.Range("D34").Value = '''0.05239'
$ws.Range("E34").Value = '''  +9.02%  '
$ws.Range("D35").Value = '''0.7840'
$ws.Range("E35").Value = '''  +13.91%  '
$ws.Range("D36").Value = '''1.221'
$ws.Range("E36").Value = '''  +8.60%  '
$ws.Range("D37").Value = '''2.806'
$ws.Range("E37").Value = '''  +3.33%  '
$ws.Range("D38").Value = '''0.02009'
$ws.Range("E38").Value = '''  +5.95%  '
$ws.Range("D39").Value = '''2.938'
$ws.Range("E39").Value = '''  +3.42%  '
$ws.Range("D40").Value = '''6.716'
$ws.Range("E40").Value = '''  +7.14%  '
$ws.Range("D41").Value = '''80.72'
$ws.Range("E41").Value = '''  +6.61%  '
$ws.Range("D42").Value = '''0.4762'
$ws.Range("E42").Value = '''  +12.85%  '
$ws.Range("D43").Value = '''2.158'
$ws.Range("E43").Value = '''  +10.43%  '
$ws.Range("D44").Value = '''107.71'
$ws.Range("E44").Value = '''  +6.92%  '
$ws.Range("D45").Value = '''0.8581'
$ws.Range("E45").Value = '''  +4.25%  '
$ws.Range("D46").Value = '''7.859'
$ws.Range("E46").Value = '''  +12.20%  '
$ws.Range("E47").Value = '''  +0.03%  '
$ws.Range("D48").Value = '''10.03'
$ws.Range("E48").Value = '''  +2.37%  '
$ws.Range("D49").Value = '''0.4367'
$ws.Range("E49").Value = '''  +11.90%  '
$ws.Range("D50").Value = '''36.94'
$ws.Range("E50").Value = '''  +5.52%  '
$ws.Range("D51").Value = '''0.1198'
$ws.Range("E51").Value = '''  +14.57%  '
